# Calendar + guest house booking data integration.
# Populates Sheet1 with the booking table (headers + 2 sample bookings),
# applies left alignment to the used range, a date number format to the
# date columns, sizes the columns to fit their content, and leaves the
# selection on the cell the author ended up on.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlLeft = -4131

# ---- Header row ---------------------------------------------------------
$headers = @(
    "Booking ID", "Date", "Customer Name", "Email", "Phone", "Guests",
    "Total Price", "Status", "Booking Date", "Special Requests",
    "Plan Name", "Plan Price"
)
for ($c = 1; $c -le $headers.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# ---- Data rows ------------------------------------------------------------
# Booking ID, Date, Customer Name, Email, Phone, Guests, Total Price,
# Status, Booking Date, Special Requests, Plan Name, Plan Price
$row2 = @("SNOW-001", 45371, "John Smith", "john@email.com", "555-0101", 2, 12800, "Confirmed", 45352, "Late check-in requested", "Weekend Getaway", 12800)
$row3 = @("SNOW-002", 45372, "John Smith", "john@email.com", "555-0101", 2, 12800, "Confirmed", 45352, "Late check-in requested", "Weekend Getaway", 12800)

for ($c = 1; $c -le $row2.Length; $c++) {
    $ws.Cells.Item(2, $c).Value = $row2[$c - 1]
}
for ($c = 1; $c -le $row3.Length; $c++) {
    $ws.Cells.Item(3, $c).Value = $row3[$c - 1]
}

# ---- Formatting -----------------------------------------------------------
# Left-align everything first (style index 1), THEN apply the date number
# format to the date columns so those cells pick up the built-in date
# format (numFmtId 14) layered on top of the left alignment (style index 2).
$used = $ws.Range("A1:L3")
$used.HorizontalAlignment = $xlLeft

$dateCells = $ws.Range("B2:B3,I2:I3")
$ws.Range("B2").NumberFormat = "mm-dd-yy"
$ws.Range("B3").NumberFormat = "mm-dd-yy"
$ws.Range("I2").NumberFormat = "mm-dd-yy"
$ws.Range("I3").NumberFormat = "mm-dd-yy"

# ---- Column widths (best-fit approximations) -------------------------------
$ws.Columns.Item(1).ColumnWidth = 9.333333333333334
$ws.Columns.Item(2).ColumnWidth = 8.666666666666666
$ws.Columns.Item(3).ColumnWidth = 14.166666666666666
$ws.Columns.Item(4).ColumnWidth = 15.0
$ws.Columns.Item(5).ColumnWidth = 7.833333333333333
$ws.Columns.Item(6).ColumnWidth = 6.0
$ws.Columns.Item(7).ColumnWidth = 9.166666666666666
$ws.Columns.Item(8).ColumnWidth = 9.333333333333334
$ws.Columns.Item(9).ColumnWidth = 11.5
$ws.Columns.Item(10).ColumnWidth = 21.0
$ws.Columns.Item(11).ColumnWidth = 16.666666666666668
$ws.Columns.Item(12).ColumnWidth = 8.666666666666666

# ---- Selection --------------------------------------------------------------
$ws.Range("E16").Select() | Out-Null
